# Added notes on excel formulas
#
# 1) Rename the existing sheet to "RGB Color Map" and add a new sheet
#    "Excel Function  Examples" (note: two spaces, matches source) that
#    becomes the active tab.
# 2) Populate the new sheet with a set of formula-examples + explanatory
#    notes in column A / column B.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename + keep as first tab, update its selection -------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "RGB Color Map"

# --- Sheet 2: brand-new "Excel Function  Examples" sheet, inserted after ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Excel Function  Examples"

# Fill content in the same order the original author would have (keeps the
# shared-string table ordering stable): column B note first, then column A
# value/formula, row by row.

# Row 1 / 2 / 3 -- INDIRECT example
$ws2.Range("B1").Value = "Cell A1 contains the name of another cell: A3"
$ws2.Range("A1").Value = "A3"
$ws2.Range("A3").Value = "Hi"
$ws2.Range("A2").Formula = "=INDIRECT(A1)"
$ws2.Range("B2").Value = "Cell A2 contains the formula: =INDIRECT(A1) which causes the value of the cell to be what is in A3"
$ws2.Range("B3").Value = "Cell A3 contains: Hi"

# Rows 5-9 -- ADDRESS examples
$ws2.Range("A5").Formula = "=ADDRESS(1,1)"
$ws2.Range("B5").Value = "Cell A5 contains the formula: =ADDRESS(1,1)"

$ws2.Range("A6").Formula = "=ADDRESS(1,1,1)"
$ws2.Range("B6").Value = "Cell A6 contains the formula: =ADDRESS(1,1,1) // Absolute"

$ws2.Range("A7").Formula = "=ADDRESS(1,1,2)"
$ws2.Range("B7").Value = "Cell A7 contains the formula: =ADDRESS(1,1,2) // Absolute Row, Relative Column"

$ws2.Range("A8").Formula = "=ADDRESS(1,1,3)"
$ws2.Range("B8").Value = "Cell A8 contains the formula: =ADDRESS(1,1,2) // Relative Row, Absolute Column"

$ws2.Range("A9").Formula = "=ADDRESS(1,1,4)"
$ws2.Range("B9").Value = "Cell A9 contains the formula: =ADDRESS(1,1,4) // Relative"

# Rows 11-14 -- CELL("address", ...) examples (12:14 share one formula)
$ws2.Range("A11").Formula = "=CELL(""address"", A1)"
$ws2.Range("B11").Value = "Cell A11 contains the formula: =CELL(""address"", A1)"

$ws2.Range("A12:A14").Formula = "=CELL(""address"", A2)"
$ws2.Range("B12").Value = "Cell A11 contains the formula: =CELL(""address"", A1)"
$ws2.Range("B13").Value = "Cell A11 contains the formula: =CELL(""address"", A1)"
$ws2.Range("B14").Value = "Cell A11 contains the formula: =CELL(""address"", A1)"

# Rows 16-17 -- CELL("row"/"col", INDIRECT(...)) examples; note column B
# here shows the formula text itself (typed with a leading apostrophe so
# it is stored as literal quote-prefixed text, not evaluated).
$ws2.Range("A16").Formula = "=CELL(""row"",INDIRECT(CELL(""address"",A1)))"
$ws2.Range("B16").Value = "'=CELL(""row"",INDIRECT(CELL(""address"",A1)))"

$ws2.Range("A17").Formula = "=CELL(""col"",INDIRECT(CELL(""address"",A1)))"
$ws2.Range("B17").Value = "'=CELL(""col"",INDIRECT(CELL(""address"",A1)))"

# Row 19
$ws2.Range("B19").Formula = "=TEXT(A11,"""")"

# --- Column widths on the new sheet -----------------------------------------
$ws2.Columns.Item(2).ColumnWidth = 88.1

# --- View state: selections + zoom + active tab -----------------------------
$ws1.Range("C2").Select()

$ws2.Activate()
$excel.ActiveWindow.Zoom = 160
$ws2.Range("A19").Select()
